# "playing with heat map"
# Fill in the latitude/longitude for the MLK incident row (row 10), which
# was previously blank, using the same coordinates already present for the
# MLK row (row 3) elsewhere in the sheet, then scroll/select so the new
# values are in view.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The two cells were blank but carried the row's shared "s=2" cell style;
# resetting to Normal first means the newly-populated cells end up with
# plain/default formatting, matching the rest of the lat/long column.
$ws.Range("F10:G10").Style = "Normal"

$ws.Range("F10").Value = 37.869249000000003
$ws.Range("G10").Value = -122.25967

# Scroll the viewport down so row 7 is at the top and select the cells we
# just filled in.
$excel.ActiveWindow.ScrollRow = 7
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("F10:G10").Select()
